# v2.2.1.10 - Provide external interface to control the MP3 module directly:
# add a new "MP3 Send Command" row to the "V2 Command" sheet that sends the
# MP3-TF-16P instruction code directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V2 Command")
$ws.Activate()

# Insert a brand-new blank row above row 46 (row 46 "Play action" and
# everything below shifts down by one, i.e. old row 46 becomes row 47, ...,
# old row 72 becomes row 73).
$ws.Rows.Item(46).Insert()

# Give the new row the same look & feel as the data rows around it:
# copy the formatting of the row that is now directly below (row 47,
# the former row 46) for columns B..K, then overwrite the values.
$ws.Range("B47:G47").Copy()
$ws.Range("B46:G46").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I47:K47").Copy()
$ws.Range("I46:K46").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new "MP3 Send Command" row content.
$ws.Range("B46").Value2 = 37
$ws.Range("C46").Value2 = "MP3 Send Command"
$ws.Range("D46").Value2 = "Yes {5}"
$ws.Range("E46").Value2 = "<cmd><parm1><parm2>"
$ws.Range("F46").Value2 = "A9 9A 05 37 12 00 01 4F ED"
$ws.Range("G46").Value2 = "A9 9A 05 37 16 00 00 52 ED"
# Columns I, J, K stay blank (just formatted) for this row, matching the
# other "direct command" rows such as the former row 48/49.

# Restore the selection to roughly where the author left it after the edit.
$ws.Range("G46").Select()
